$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Investments")

# Update GLD max weight constraint from 5% to 1%
$ws.Range("C19").Value = 0.01

# Reflect the active cell/selection as saved in the authored workbook
$ws.Activate()
$ws.Range("C19").Select()
